$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; this shifts the existing rows 85:167
# (and everything below) down by one, to 86:168, preserving their values,
# matching the dimension growing from A1:T167 to A1:T168.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record's data.
$ws.Range("A85").Value = 10
$ws.Range("B85").Value = "Vega Modelo de Temuco"
$ws.Range("C85").Value = "La Araucanía"
$ws.Range("D85").Value = 44874
$ws.Range("D85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E85").Value = 9
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100107
$ws.Range("H85").Value = "Otros"
$ws.Range("I85").Value = 100107002
$ws.Range("J85").Value = "Chirimoya"
$ws.Range("K85").Value = "Cultivar IV Región"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 200
$ws.Range("N85").Value = 2800
$ws.Range("O85").Value = 2800
$ws.Range("P85").Value = 2800
$ws.Range("Q85").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R85").Value = "Provincia del Elquí"
$ws.Range("S85").Value = 2800
$ws.Range("T85").Value = 1
